# Refresh generated notebook-run dataset: update this round's scores in
# column B ("Pontuacoes") for the teams whose totals changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 41.36
$ws.Range("B3").Value = 57.26
$ws.Range("B4").Value = 59.69
$ws.Range("B5").Value = 49.36
$ws.Range("B6").Value = 56.09
$ws.Range("B8").Value = 43.56
$ws.Range("B9").Value = 58.26
$ws.Range("B11").Value = 71.36
$ws.Range("B12").Value = 81.76000000000001
$ws.Range("B14").Value = 54.95
$ws.Range("B16").Value = 66.86
$ws.Range("B17").Value = 54.66
$ws.Range("B18").Value = 63.76
$ws.Range("B19").Value = 47.86
$ws.Range("B20").Value = 71.16
$ws.Range("B21").Value = 50.85
